$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "comment" table: remove the "heading" field, shifting content/time up one row
$ws.Range("G12").Value = "content"
$ws.Range("G13").Value = "time"
$ws.Range("G14").ClearContents()

# "payment" table: remove the "amount" field, shifting type up one row
$ws.Range("K12").Value = "type"
$ws.Range("K13").ClearContents()

# "cart" table: swap Vietnamese labels for English equivalents
$ws.Range("E15").Value = "phone"
$ws.Range("E16").Value = "describe"

# "product" table: drop the trailing "status" field/row
$ws.Range("C18").ClearContents()

# New note cell far below the tables
$ws.Range("O23").Value = " "

# Update the active selection shown when the workbook is reopened
[void]$ws.Range("L9").Select()
